$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Cells.Item(1, 6).Value = "Last status check on: 25.02.2022 09:15"

# Update row 8 (Benzina Albert Modrice) values
$ws.Cells.Item(8, 2).Value = 37.9
$ws.Cells.Item(8, 3).Value = 37.5

# D8 becomes a literal text "+0.4" (not a number) with no special style
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "+0.4"
$ws.Cells.Item(8, 4).Style = "Normal"

# E8 becomes a literal text timestamp (not a date number) with no special style
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "2022-02-25 09:17:23"
$ws.Cells.Item(8, 5).Style = "Normal"
